# Updated cryptos list values (price + 1h volume %) per upstream refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.030.62"
$ws.Range("E2").Value = "  -3.38%  "

$ws.Range("D3").Value = "'1.843.83"
$ws.Range("E3").Value = "  -2.44%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").Value = "'0.7007"
$ws.Range("E5").Value = "  -6.02%  "

$ws.Range("D6").Value = "'237.62"
$ws.Range("E6").Value = "  -2.22%  "

$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.34%  "

$ws.Range("D8").Value = "'0.3035"
$ws.Range("E8").Value = "  -4.24%  "

$ws.Range("D9").Value = "'0.07417"
$ws.Range("E9").Value = "  +2.29%  "

$ws.Range("D10").Value = "'23.20"
$ws.Range("E10").Value = "  -7.42%  "

$ws.Range("D11").Value = "'0.08115"
$ws.Range("E11").Value = "  -2.93%  "

$ws.Range("D12").Value = "'0.7228"
$ws.Range("E12").Value = "  -4.92%  "

$ws.Range("D13").Value = "'5.217"
$ws.Range("E13").Value = "  -3.88%  "

$ws.Range("D14").Value = "'1.813.27"
$ws.Range("E14").Value = "  -3.65%  "

$ws.Range("D15").Value = "'88.81"
$ws.Range("E15").Value = "  -4.19%  "

$ws.Range("D16").Value = "'28.930.79"
$ws.Range("E16").Value = "  -3.52%  "

$ws.Range("D17").Value = "'5.778"
$ws.Range("E17").Value = "  -6.18%  "

$ws.Range("D18").Value = "'239.34"
$ws.Range("E18").Value = "  -4.21%  "

$ws.Range("D19").Value = "'0.000007655"
$ws.Range("E19").Value = "  -2.85%  "

$ws.Range("D20").Value = "'12.99"
$ws.Range("E20").Value = "  -4.87%  "

$ws.Range("D21").Value = "'0.9996"
$ws.Range("E21").Value = "  +0.15%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.29%  "

$ws.Range("D23").Value = "'2.056.25"
$ws.Range("E23").Value = "  -1.95%  "

$ws.Range("D24").Value = "'7.516"
$ws.Range("E24").Value = "  -6.13%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'160.91"
$ws.Range("E25").Value = "  -2.73%  "

$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").Value = "'0.1455"
$ws.Range("E26").Value = "  -8.26%  "

$ws.Range("D27").Value = "'8.930"
$ws.Range("E27").Value = "  -4.12%  "

$ws.Range("D28").Value = "'17.94"
$ws.Range("E28").Value = "  -4.41%  "

$ws.Range("D29").Value = "'1.929"
$ws.Range("E29").Value = "  -5.91%  "

$ws.Range("D30").Value = "'1.376"
$ws.Range("E30").Value = "  -7.60%  "

$ws.Range("D31").Value = "'4.481"
$ws.Range("E31").Value = "  -2.94%  "

$ws.Range("E32").Value = "  -2.93%  "

$ws.Range("D33").Value = "'4.003"
$ws.Range("E33").Value = "  -5.46%  "

$ws.Range("D34").Value = "'0.05148"
$ws.Range("E34").Value = "  -4.39%  "

$ws.Range("D35").Value = "'1.181"
$ws.Range("E35").Value = "  -5.90%  "

$ws.Range("D36").Value = "'0.7036"
$ws.Range("E36").Value = "  -7.36%  "

$ws.Range("E37").Value = "  +3.76%  "

$ws.Range("D38").Value = "'2.643"
$ws.Range("E38").Value = "  -2.42%  "

$ws.Range("D39").Value = "'0.01862"
$ws.Range("E39").Value = "  -5.41%  "

$ws.Range("D40").Value = "'2.674"
$ws.Range("E40").Value = "  -3.37%  "

$ws.Range("D41").Value = "'0.9002"
$ws.Range("E41").Value = "  +3.28%  "

$ws.Range("D42").Value = "'5.961"
$ws.Range("E42").Value = "  -1.28%  "

$ws.Range("D43").Value = "'0.4268"
$ws.Range("E43").Value = "  -6.80%  "

$ws.Range("D44").Value = "'1.059.11"
$ws.Range("E44").Value = "  -4.10%  "

$ws.Range("D45").Value = "'69.70"
$ws.Range("E45").Value = "  -4.44%  "

$ws.Range("D46").Value = "'0.9998"
$ws.Range("E46").Value = "  +0.07%  "

$ws.Range("D47").Value = "'101.69"
$ws.Range("E47").Value = "  -2.78%  "

$ws.Range("D48").Value = "'1.746"
$ws.Range("E48").Value = "  -6.76%  "

$ws.Range("E49").Value = "  -4.58%  "

$ws.Range("D50").Value = "'7.025"
$ws.Range("E50").Value = "  -7.93%  "

$ws.Range("D51").Value = "'1.967.98"
$ws.Range("E51").Value = "  -3.30%  "
